$d = $word.ActiveDocument

# Locate the paragraph that begins "Ver no Jupiter Salvar em pdf Salvar em docx".
$r1 = $d.Content
$found1 = $r1.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startIndex = $r1.Paragraphs.First.Index

# Locate the paragraph that begins with the "(c) 2020 ..." copyright notice.
$r2 = $d.Content
$found2 = $r2.Find.Execute([char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endIndex = $r2.Paragraphs.First.Index

if ($found1 -and $found2) {
    # Also remove the blank paragraph immediately preceding "Ver no Jupiter...",
    # so the entire footer block (blank line, "Ver no Jupiter..." line, and the
    # copyright line) is deleted, paragraph marks included.
    $delStart = $d.Paragraphs.Item($startIndex - 1).Range.Start
    $delEnd = $d.Paragraphs.Item($endIndex).Range.End

    $d.Range($delStart, $delEnd).Delete()
}
